$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("B3").Value = 432
$ws.Range("C3").Value = 311
$ws.Range("D3").Value = 87
$ws.Range("E3").Value = 42
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 9

$ws2 = $wb.Worksheets.Item("DEF")
$ws2.Range("B3").Value = 356
$ws2.Range("C3").Value = 246
$ws2.Range("D3").Value = 101
$ws2.Range("E3").Value = 44
